$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "abfef73e-eb84-458c-bdca-097916493317.md" file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 10:15:19"

# --- zh-cn sheet: row 3 is the "abfef73e-eb84-458c-bdca-097916493317.md" file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-28 10:15:15"

# --- de-de sheet: row 3 is the "abfef73e-eb84-458c-bdca-097916493317.md" file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-28 10:15:19"

# --- Column widths widen to fit the longer "Ready for handoff" status text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
